# Normalize the two worksheet (tab) names to lowercase, matching the new
# CSV-based naming convention ("sedementation" -> "sedimentation",
# "Eutrophication" -> "eutrophication"). Renaming a sheet automatically
# updates any defined names / formulas that reference it by name.
$wb = $excel.ActiveWorkbook

$wsSed = $wb.Worksheets.Item("sedementation")
$wsSed.Name = "sedimentation"

$wsEut = $wb.Worksheets.Item("Eutrophication")
$wsEut.Name = "eutrophication"

# Switch the active/selected tab from "Eutrophication" to "sedimentation"
# (the first sheet).
$wsSed.Activate()

$wb.Save()
